$d = $word.ActiveDocument

# The cover letter uses the hyphenated form "cross-sections" three times
# ("...puttering cross-sections in 2D crystals...", "...sputtering
# cross-sections in 2D crystals by explicitly calculating...", and
# "...yields cross-sections that quantitatively...").  The commit removes
# the dash in all three spots, turning "cross-sections" into "cross
# sections".
#
# We replace only the single hyphen character (not the whole phrase) and
# re-search from scratch after each substitution.  This is the most
# surgical edit possible (touches exactly one character per occurrence)
# and keeps every other run in the paragraph completely untouched.

$replacedCount = 0
while ($true) {
    $text = $d.Content.Text
    $idx = $text.IndexOf("cross-sections")
    if ($idx -lt 0) { break }

    $hyphenPos = $idx + 5   # "cross" is 5 characters, so the hyphen sits right after it
    $rng = $d.Range($hyphenPos, $hyphenPos + 1)
    $rng.Text = " "

    $replacedCount = $replacedCount + 1
    if ($replacedCount -gt 50) { break }   # safety valve against infinite loops
}

Write-Host "Removed the dash in $replacedCount occurrence(s) of 'cross-sections'."
